$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 29; this shifts existing rows 29..110 down to 30..111,
# which reproduces the data exactly (each old row's data moves one row down,
# and the previous last row 110 becomes row 111 unchanged).
$ws.Rows("29:29").Insert()

# Populate the new row 29 with the new weekly data point.
$ws.Range("A29").Value = 11
$ws.Range("B29").Value = "Vega Monumental Concepción"
$ws.Range("C29").Value = "Bíobío"
$ws.Range("D29").Value = 44497
$ws.Range("E29").Value = 8
$ws.Range("F29").Value = "Fruta"
$ws.Range("G29").Value = 100108
$ws.Range("H29").Value = "Tropicales y subtropicales"
$ws.Range("I29").Value = 100108005
$ws.Range("J29").Value = "Piña"
$ws.Range("K29").Value = "Caramelo"
$ws.Range("L29").Value = "Primera"
$ws.Range("M29").Value = 200
$ws.Range("N29").Value = 19000
$ws.Range("O29").Value = 20000
$ws.Range("P29").Value = 19500
$ws.Range("Q29").Value = "$/caja 12 unidades"
$ws.Range("R29").Value = "Ecuador"
$ws.Range("S29").Value = 1625
$ws.Range("T29").Value = 12
